# Generate Report for Archive
# The localization status for the two handed-off files moved on from
# "Ready for handoff" to "In Translation" - update every occurrence across
# the Overview sheet (zh-cn/de-de status columns) and each per-language
# sheet's Status column. Excel then re-flows (narrows) those status
# columns to fit the new, shorter text.

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Overview sheet: columns E (zh-cn) and F (de-de) hold the per-language
# status for each of the two files (rows 2 and 3).
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Per-language sheets: column C is "Status" for each of the two files.
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# The status text got shorter ("Ready for handoff" -> "In Translation"),
# so the status columns shrink to fit the new content.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
